$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data cells in C2:F8 hold numbers that are stored as text (as in the
# original workbook, where every value cell is typed t="str"). Force a
# text number format so Excel keeps the new values as text too.
$rng = $ws.Range("C2:F8")
$rng.NumberFormat = "@"

# New values for runs/balls/fours/sixes (columns C,D,E,F) for rows 2-8,
# per the target diff (rows 2-8 values get permuted).
$values = @(
    @("11", "7", "1", "1"),  # row 2
    @("1",  "2", "0", "0"),  # row 3
    @("0",  "0", "0", "0"),  # row 4
    @("14", "8", "1", "1"),  # row 5
    @("0",  "1", "0", "0"),  # row 6
    @("6",  "5", "1", "0"),  # row 7
    @("3",  "7", "0", "0")   # row 8
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 4; $j++) {
        $col = $j + 3  # column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$i][$j]
    }
}
